$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster was re-sorted/re-ranked (rows 8-12, 15-17 reshuffled) and a new
# player (Omer Yurtseven, previously row 16) now sits in a fresh row 18,
# with Cody Zeller added as the new row 15.

# Row 8: now Haywood Highsmith (was Kyle Lowry)
$ws.Cells.Item(8, 2).Value = 24
$ws.Cells.Item(8, 3).Value = "Haywood Highsmith"
$ws.Cells.Item(8, 4).Value = "SF"
$ws.Cells.Item(8, 5).Value = "6-7"
$ws.Cells.Item(8, 6).Value = 220
$ws.Cells.Item(8, 7).Value = "December 9, 1996"
$ws.Cells.Item(8, 8).Value = "us"
$ws.Cells.Item(8, 9).Value = "2"
$ws.Cells.Item(8, 10).Value = "Wheeling University"
$ws.Cells.Item(8, 11).Value = "https://www.basketball-reference.com/players/h/highsha01.html"

# Row 9: now Kyle Lowry (was Haywood Highsmith)
$ws.Cells.Item(9, 2).Value = 7
$ws.Cells.Item(9, 3).Value = "Kyle Lowry"
$ws.Cells.Item(9, 4).Value = "PG"
$ws.Cells.Item(9, 5).Value = "6-0"
$ws.Cells.Item(9, 6).Value = 196
$ws.Cells.Item(9, 7).Value = "March 25, 1986"
$ws.Cells.Item(9, 8).Value = "us"
$ws.Cells.Item(9, 9).Value = "16"
$ws.Cells.Item(9, 10).Value = "Villanova"
$ws.Cells.Item(9, 11).Value = "https://www.basketball-reference.com/players/l/lowryky01.html"

# Row 10: now Victor Oladipo (was Duncan Robinson)
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 3).Value = "Victor Oladipo"
$ws.Cells.Item(10, 4).Value = "SG"
$ws.Cells.Item(10, 5).Value = "6-4"
$ws.Cells.Item(10, 6).Value = 213
$ws.Cells.Item(10, 7).Value = "May 4, 1992"
$ws.Cells.Item(10, 8).Value = "us"
$ws.Cells.Item(10, 9).Value = "9"
$ws.Cells.Item(10, 10).Value = "Indiana"
$ws.Cells.Item(10, 11).Value = "https://www.basketball-reference.com/players/o/oladivi01.html"

# Row 11: now Duncan Robinson (was Orlando Robinson (TW))
$ws.Cells.Item(11, 2).Value = 55
$ws.Cells.Item(11, 3).Value = "Duncan Robinson"
$ws.Cells.Item(11, 4).Value = "SF"
$ws.Cells.Item(11, 5).Value = "6-7"
$ws.Cells.Item(11, 6).Value = 215
$ws.Cells.Item(11, 7).Value = "April 22, 1994"
$ws.Cells.Item(11, 8).Value = "us"
$ws.Cells.Item(11, 9).Value = "4"
$ws.Cells.Item(11, 10).Value = "Williams, Michigan"
$ws.Cells.Item(11, 11).Value = "https://www.basketball-reference.com/players/r/robindu01.html"

# Row 12: now Orlando Robinson (TW) (was Victor Oladipo)
$ws.Cells.Item(12, 2).Value = 25
$ws.Cells.Item(12, 3).Value = "Orlando Robinson (TW)"
$ws.Cells.Item(12, 4).Value = "C"
$ws.Cells.Item(12, 5).Value = "7-0"
$ws.Cells.Item(12, 6).Value = 235
$ws.Cells.Item(12, 7).Value = "July 10, 2000"
$ws.Cells.Item(12, 8).Value = "us"
$ws.Cells.Item(12, 9).Value = "R"
$ws.Cells.Item(12, 10).Value = "Fresno State"
$ws.Cells.Item(12, 11).Value = "https://www.basketball-reference.com/players/r/robinor01.html"

# Rows 13 (Nikola Jović) and 14 (Jamal Cain (TW)) are unchanged.

# Row 15: now Cody Zeller (was Udonis Haslem)
$ws.Cells.Item(15, 2).Value = 44
$ws.Cells.Item(15, 3).Value = "Cody Zeller"
$ws.Cells.Item(15, 4).Value = "C"
$ws.Cells.Item(15, 5).Value = "6-11"
$ws.Cells.Item(15, 6).Value = 240
$ws.Cells.Item(15, 7).Value = "October 5, 1992"
$ws.Cells.Item(15, 8).Value = "us"
$ws.Cells.Item(15, 9).Value = "9"
$ws.Cells.Item(15, 10).Value = "Indiana"
$ws.Cells.Item(15, 11).Value = "https://www.basketball-reference.com/players/z/zelleco01.html"

# Row 16: now Kevin Love (was Omer Yurtseven)
$ws.Cells.Item(16, 2).Value = 42
$ws.Cells.Item(16, 3).Value = "Kevin Love"
$ws.Cells.Item(16, 4).Value = "PF"
$ws.Cells.Item(16, 5).Value = "6-8"
$ws.Cells.Item(16, 6).Value = 251
$ws.Cells.Item(16, 7).Value = "September 7, 1988"
$ws.Cells.Item(16, 8).Value = "us"
$ws.Cells.Item(16, 9).Value = "14"
$ws.Cells.Item(16, 10).Value = "UCLA"
$ws.Cells.Item(16, 11).Value = "https://www.basketball-reference.com/players/l/loveke01.html"

# Row 17: now Udonis Haslem (was Kevin Love)
$ws.Cells.Item(17, 2).Value = 40
$ws.Cells.Item(17, 3).Value = "Udonis Haslem"
$ws.Cells.Item(17, 4).Value = "C"
$ws.Cells.Item(17, 5).Value = "6-8"
$ws.Cells.Item(17, 6).Value = 235
$ws.Cells.Item(17, 7).Value = "June 9, 1980"
$ws.Cells.Item(17, 8).Value = "us"
$ws.Cells.Item(17, 9).Value = "19"
$ws.Cells.Item(17, 10).Value = "Florida"
$ws.Cells.Item(17, 11).Value = "https://www.basketball-reference.com/players/h/hasleud01.html"

# Row 18 is brand new: Omer Yurtseven (was row 16)
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 1).Style = $ws.Cells.Item(17, 1).Style
$ws.Cells.Item(18, 3).Value = "Omer Yurtseven"
$ws.Cells.Item(18, 4).Value = "C"
$ws.Cells.Item(18, 5).Value = "7-0"
$ws.Cells.Item(18, 6).Value = 264
$ws.Cells.Item(18, 7).Value = "June 19, 1998"
$ws.Cells.Item(18, 8).Value = "tr"
$ws.Cells.Item(18, 9).Value = "1"
$ws.Cells.Item(18, 10).Value = "NC State, Georgetown"
$ws.Hyperlinks.Add($ws.Cells.Item(18, 11), "https://www.basketball-reference.com/players/y/yurtsom01.html")
$ws.Cells.Item(18, 11).Style = $ws.Cells.Item(17, 11).Style
